$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Constants for alignment (Excel COM constants)
# ---------------------------------------------------------------------------
$xlCenter = -4108
$xlRight  = -4152

# ---------------------------------------------------------------------------
# RAW sheet (sheet1.xml) - add new funcionarios_2023_01 dataset summary
# ---------------------------------------------------------------------------
$raw = $wb.Worksheets.Item("RAW")

# New header cells in row 1 (bold, centered - same visual style as existing A1/B1)
$raw.Range("B1").Value = "BRUTO"

$raw.Range("C1").Font.Bold = $true
$raw.Range("C1").HorizontalAlignment = $xlCenter
$raw.Range("C1").VerticalAlignment = $xlCenter
$raw.Range("C1").Value = "ELIMINADO"

$raw.Range("D1").Font.Bold = $true
$raw.Range("D1").HorizontalAlignment = $xlCenter
$raw.Range("D1").VerticalAlignment = $xlCenter
$raw.Range("D1").Value = "CRUDO"

# Dataset name cell A2 (plain vertical-centered style)
$raw.Range("A2").VerticalAlignment = $xlCenter
$raw.Range("A2").Value = "funcionarios_2023_01"

# Numeric cells with thousands separator format (#,##0), vertical centered
$numCells = @("B2","C2","I2","I3","I4","I5","I6","I7","I8","I9","I10","I11","I13","I14")
foreach ($addr in $numCells) {
    $raw.Range($addr).NumberFormat = "#,##0"
    $raw.Range($addr).VerticalAlignment = $xlCenter
}

$raw.Range("B2").Value = 862950
$raw.Range("C2").Value = 6825

$raw.Range("I2").Value = 148522
$raw.Range("I3").Value = 166164
$raw.Range("I4").Value = 204020
$raw.Range("I5").Value = 120775
$raw.Range("I6").Value = 40029
$raw.Range("I7").Value = 17699
$raw.Range("I8").Value = 13442
$raw.Range("I9").Value = 14448
$raw.Range("I10").Value = 15684
$raw.Range("I11").Value = 120
$raw.Range("I13").Value = 8498
$raw.Range("I14").Value = 106724

# Formulas (also #,##0 + vertical centered, same style family as numeric cells)
$formulaCells = @("D2","J10","J15")
foreach ($addr in $formulaCells) {
    $raw.Range($addr).NumberFormat = "#,##0"
    $raw.Range($addr).VerticalAlignment = $xlCenter
}
$raw.Range("D2").Formula = "=B2-C2"
$raw.Range("J10").Formula = "=+SUM(I2:I10)"
$raw.Range("J15").Formula = "=+SUM(I13:I14)"

# Plain vertical-centered numeric cell (J11) - no thousands separator
$raw.Range("J11").VerticalAlignment = $xlCenter
$raw.Range("J11").Value = 120

# Plain vertical-centered text cells (K10, K11)
$raw.Range("K10").VerticalAlignment = $xlCenter
$raw.Range("K10").Value = "PY"
$raw.Range("K11").VerticalAlignment = $xlCenter
$raw.Range("K11").Value = "EXT"

# Right-aligned / vertical-centered index + label cells in column H
$hCells = @("H2","H3","H4","H5","H6","H7","H8","H9","H10")
$hVals  = @(1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $hCells.Count; $i++) {
    $raw.Range($hCells[$i]).HorizontalAlignment = $xlRight
    $raw.Range($hCells[$i]).VerticalAlignment = $xlCenter
    $raw.Range($hCells[$i]).Value = $hVals[$i]
}

$raw.Range("H11").HorizontalAlignment = $xlRight
$raw.Range("H11").VerticalAlignment = $xlCenter
$raw.Range("H11").Value = "E"

$raw.Range("H13").HorizontalAlignment = $xlRight
$raw.Range("H13").VerticalAlignment = $xlCenter
$raw.Range("H13").Value = "V"

$raw.Range("H14").HorizontalAlignment = $xlRight
$raw.Range("H14").VerticalAlignment = $xlCenter
$raw.Range("H14").Value = "A"

# Column A width grows to fit the new longer dataset name
$raw.Columns(1).ColumnWidth = 17.830729166666668

# Selection / active sheet
$raw.Range("D7").Select()
$raw.Activate()

Write-Host "RAW sheet updated"

# ---------------------------------------------------------------------------
# OLTP sheet (sheet3.xml) - no longer the active tab, selection stays I15
# ---------------------------------------------------------------------------
$oltp = $wb.Worksheets.Item("OLTP")
$oltp.Range("I15").Select()

Write-Host "done"
